$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'32.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.48%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.074"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.00%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07714"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.56%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.022"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-10.39%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.899"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.805"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.12%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9264"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.06%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-0.28%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08194"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'4.90%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08615"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.02%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03090"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.29%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1000"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.32%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.62%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005900"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.85%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'3.477"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.54%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-4.04%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'1.88%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.09%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.396"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.23%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'9.95%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04544"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.23%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-1.81%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.88%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001254"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.07%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01734"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.94%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.86%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007490"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.87%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1360"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-1.00%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002337"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.48%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01042"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'4.74%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'-1.45%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.11%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.800"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'127.79%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-16.79%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.11%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("E50").Style = "Normal"
